$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("B1").Value2 = "Categoría (mm)"
$ws.Range("C1").Value2 = "Importancia"
$ws.Range("D1").Value2 = "FV"

# --- Data rows: re-labelled categories and re-ordered IDs ---
# Row 2 (Muy alta)
$ws.Range("A2").Value2 = 5
$ws.Range("B2").Value2 = "1081-1233"

# Row 3 (Alta)
$ws.Range("A3").Value2 = 4
$ws.Range("B3").Value2 = "930-1081"

# Row 4 (Moderada)
$ws.Range("A4").Value2 = 3
$ws.Range("B4").Value2 = "778-930"

# Row 5 (Baja)
$ws.Range("A5").Value2 = 2
$ws.Range("B5").Value2 = "626-778"

# Row 6 (Baja)
$ws.Range("A6").Value2 = 1
$ws.Range("B6").Value2 = "551-626"

# --- Remove the extra blank template row 7 ---
$ws.Range("A7:D7").Clear()

# --- Selection ---
$ws.Range("B8").Select()

$wb.Save()
